# Better validation for aviation: the INCLUDES_HULL (AM) / INCLUDES_LIABILITY (AN)
# flags on the "sections" sheet were stored as TRUE/FALSE booleans. Re-author them
# as plain 0/1 numeric values instead, and drop the flags entirely for the rows
# that shouldn't carry them (rows 46-98).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sections")

# Exceptions where AM/AN aren't both 1 (row -> [AM, AN])
$exceptions = @{
    38 = @(1, 0)
    39 = @(0, 1)
}

# Rows 2-45: re-write AM/AN as numeric 0/1 (was boolean TRUE/FALSE) keeping values.
for ($r = 2; $r -le 45; $r++) {
    if ($exceptions.ContainsKey($r)) {
        $pair = $exceptions[$r]
        $amVal = $pair[0]
        $anVal = $pair[1]
    } else {
        $amVal = 1
        $anVal = 1
    }
    $ws.Cells.Item($r, 39).Value = $amVal   # column AM
    $ws.Cells.Item($r, 40).Value = $anVal   # column AN
}

# Rows 46-98: these rows shouldn't carry the AM/AN flags at all - clear them.
for ($r = 46; $r -le 98; $r++) {
    $ws.Range("AM" + $r + ":AN" + $r).ClearContents()
}
